$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Opcodes")
$ws.Select()
$sv = $excel.ActiveWindow.ActiveSheetView
Write-Host ($sv | Get-Member | Out-String)
